$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the blank separator row above the "SEZON" price table (old row 10).
# This shifts every subsequent row up by one, matching the target layout.
$ws.Rows.Item(10).Delete()

# The "SINGLE" room-type label (now on row 11) becomes "SINGLE ROOM".
$ws.Range("A11").Value = "SINGLE ROOM"

# The print area shrank by one row along with the deleted row.
$ws.PageSetup.PrintArea = '$A$1:$H$37'

# Restore the active selection shown in the edited workbook.
$ws.Range("A11").Select()
